# Auto-generated Excel COM-interop script to apply the market-data refresh diff
# to the Coeurl_Profits workbook. Operates per-sheet, per-cell: value updates,
# one cell addition (N58 on CUL), and a few cell removals (ClearContents) where
# the target row no longer carries a LeveProfit value for that column.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 4618   # H9: 5124.9443 -> 4618
$ws.Cells.Item(9, 9).Value = 5746.25   # I9: 6126.7334 -> 5746.25
$ws.Cells.Item(9, 10).Value = 105   # J9: 116 -> 105
$ws.Cells.Item(9, 11).Value = 5746.25   # K9: 6126.7334 -> 5746.25
$ws.Cells.Item(9, 12).Value = 105   # L9: 116 -> 105
$ws.Cells.Item(9, 13).Value = -5577.25   # M9: -5957.7334 -> -5577.25
$ws.Cells.Item(9, 14).Value = -443   # N9: -454 -> -443
$ws.Cells.Item(19, 8).Value = 1291.125   # H19: 1267.9048 -> 1291.125
$ws.Cells.Item(19, 9).Value = 995.6667   # I19: 929.8 -> 995.6667
$ws.Cells.Item(19, 10).Value = 1359.3077   # J19: 1373.5625 -> 1359.3077
$ws.Cells.Item(19, 11).Value = 995.6667   # K19: 929.8 -> 995.6667
$ws.Cells.Item(19, 12).Value = 1359.3077   # L19: 1373.5625 -> 1359.3077
$ws.Cells.Item(19, 13).Value = -820.6667   # M19: -754.8 -> -820.6667
$ws.Cells.Item(19, 14).Value = -1709.3077   # N19: -1723.5625 -> -1709.3077
$ws.Cells.Item(112, 8).Value = 92562.73   # H112: 101669.5 -> 92562.73
$ws.Cells.Item(112, 10).Value = 101709   # J112: 112843.89 -> 101709
$ws.Cells.Item(112, 12).Value = 305127   # L112: 338531.67 -> 305127
$ws.Cells.Item(112, 14).Value = -307343   # N112: -340747.67 -> -307343
$ws.Cells.Item(125, 8).Value = 2333   # H125: 2220 -> 2333
$ws.Cells.Item(125, 10).Value = 2999.5   # J125: 2626.6667 -> 2999.5
$ws.Cells.Item(125, 12).Value = 26995.5   # L125: 23640.0003 -> 26995.5
$ws.Cells.Item(125, 14).Value = -31915.5   # N125: -28560.0003 -> -31915.5
$ws.Cells.Item(138, 8).Value = 4174.3184   # H138: 4327.2383 -> 4174.3184
$ws.Cells.Item(138, 9).Value = 1452.6428   # I138: 1452.7142 -> 1452.6428
$ws.Cells.Item(138, 10).Value = 5444.433   # J138: 5764.5 -> 5444.433
$ws.Cells.Item(138, 11).Value = 4357.928400000001   # K138: 4358.142599999999 -> 4357.928400000001
$ws.Cells.Item(138, 12).Value = 16333.299   # L138: 17293.5 -> 16333.299
$ws.Cells.Item(138, 13).Value = 782.0715999999993   # M138: 781.8574000000008 -> 782.0715999999993
$ws.Cells.Item(138, 14).Value = -26613.299   # N138: -27573.5 -> -26613.299
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 45729.684   # H5: 50346.65 -> 45729.684
$ws.Cells.Item(5, 9).Value = 77077.30499999999   # I5: 91171.37 -> 77077.30499999999
$ws.Cells.Item(5, 11).Value = 77077.30499999999   # K5: 91171.37 -> 77077.30499999999
$ws.Cells.Item(5, 13).Value = -76965.30499999999   # M5: -91059.37 -> -76965.30499999999
$ws.Cells.Item(32, 8).Value = 5401.7744   # H32: 5716.276 -> 5401.7744
$ws.Cells.Item(32, 9).Value = 4931.6665   # I32: 5258.88 -> 4931.6665
$ws.Cells.Item(32, 11).Value = 4931.6665   # K32: 5258.88 -> 4931.6665
$ws.Cells.Item(32, 13).Value = -4644.6665   # M32: -4971.88 -> -4644.6665
$ws.Cells.Item(56, 8).Value = 7083.3335   # H56: 17500 -> 7083.3335
$ws.Cells.Item(56, 9).Value = 6304.3477   # I56: 15000 -> 6304.3477
$ws.Cells.Item(56, 11).Value = 6304.3477   # K56: 15000 -> 6304.3477
$ws.Cells.Item(56, 13).Value = -5562.3477   # M56: -14258 -> -5562.3477
$ws.Cells.Item(88, 8).Value = 1998.6364   # H88: 1844.5 -> 1998.6364
$ws.Cells.Item(88, 9).Value = 3117.2   # I88: 2622.5 -> 3117.2
$ws.Cells.Item(88, 11).Value = 3117.2   # K88: 2622.5 -> 3117.2
$ws.Cells.Item(88, 13).Value = -2711.2   # M88: -2216.5 -> -2711.2
$ws.Cells.Item(91, 8).Value = 1998.6364   # H91: 1844.5 -> 1998.6364
$ws.Cells.Item(91, 9).Value = 3117.2   # I91: 2622.5 -> 3117.2
$ws.Cells.Item(91, 11).Value = 3117.2   # K91: 2622.5 -> 3117.2
$ws.Cells.Item(91, 13).Value = -1713.2   # M91: -1218.5 -> -1713.2
$ws.Cells.Item(97, 8).Value = 441.66666   # H97: 824.8611 -> 441.66666
$ws.Cells.Item(97, 9).Value = 0   # I97: 952 -> 0
$ws.Cells.Item(97, 10).Value = 441.66666   # J97: 298.14285 -> 441.66666
$ws.Cells.Item(97, 11).Value = 0   # K97: 952 -> 0
$ws.Cells.Item(97, 12).Value = 441.66666   # L97: 298.14285 -> 441.66666
$ws.Cells.Item(97, 13).ClearContents()   # M97: remove (was -456)
$ws.Cells.Item(97, 14).Value = -1433.66666   # N97: -1290.14285 -> -1433.66666
$ws.Cells.Item(132, 8).Value = 2311.4888   # H132: 2493.195 -> 2311.4888
$ws.Cells.Item(132, 9).Value = 1861.1842   # I132: 2027.2941 -> 1861.1842
$ws.Cells.Item(132, 10).Value = 4756   # J132: 4756.143 -> 4756
$ws.Cells.Item(132, 11).Value = 5583.5526   # K132: 6081.8823 -> 5583.5526
$ws.Cells.Item(132, 12).Value = 14268   # L132: 14268.429 -> 14268
$ws.Cells.Item(132, 13).Value = -3053.5526   # M132: -3551.8823 -> -3053.5526
$ws.Cells.Item(132, 14).Value = -19328   # N132: -19328.429 -> -19328
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 45729.684   # H4: 50346.65 -> 45729.684
$ws.Cells.Item(4, 9).Value = 77077.30499999999   # I4: 91171.37 -> 77077.30499999999
$ws.Cells.Item(4, 11).Value = 77077.30499999999   # K4: 91171.37 -> 77077.30499999999
$ws.Cells.Item(4, 13).Value = -76962.30499999999   # M4: -91056.37 -> -76962.30499999999
$ws.Cells.Item(22, 8).Value = 217137.6   # H22: 224619.89 -> 217137.6
$ws.Cells.Item(22, 9).Value = 341   # I22: 358.27274 -> 341
$ws.Cells.Item(22, 11).Value = 341   # K22: 358.27274 -> 341
$ws.Cells.Item(22, 13).Value = -168   # M22: -185.27274 -> -168
$ws.Cells.Item(134, 8).Value = 2006.36   # H134: 2001.9387 -> 2006.36
$ws.Cells.Item(134, 9).Value = 1783   # I134: 1784.674 -> 1783
$ws.Cells.Item(134, 10).Value = 4575   # J134: 5333.3335 -> 4575
$ws.Cells.Item(134, 11).Value = 5349   # K134: 5354.022 -> 5349
$ws.Cells.Item(134, 12).Value = 13725   # L134: 16000.0005 -> 13725
$ws.Cells.Item(134, 13).Value = -2814   # M134: -2819.022 -> -2814
$ws.Cells.Item(134, 14).Value = -18795   # N134: -21070.0005 -> -18795
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 8244.666999999999   # H28: 9033.333000000001 -> 8244.666999999999
$ws.Cells.Item(28, 10).Value = 7900.375   # J28: 8787.625 -> 7900.375
$ws.Cells.Item(28, 12).Value = 7900.375   # L28: 8787.625 -> 7900.375
$ws.Cells.Item(28, 14).Value = -8390.375   # N28: -9277.625 -> -8390.375
$ws.Cells.Item(99, 8).Value = 3562.3845   # H99: 3634.25 -> 3562.3845
$ws.Cells.Item(99, 9).Value = 3442.5833   # I99: 3510.0908 -> 3442.5833
$ws.Cells.Item(99, 11).Value = 3442.5833   # K99: 3510.0908 -> 3442.5833
$ws.Cells.Item(99, 13).Value = -1944.5833   # M99: -2012.0908 -> -1944.5833
$ws.Cells.Item(115, 8).Value = 40999.332   # H115: 68499.164 -> 40999.332
$ws.Cells.Item(115, 10).Value = 40999.332   # J115: 68499.164 -> 40999.332
$ws.Cells.Item(115, 12).Value = 40999.332   # L115: 68499.164 -> 40999.332
$ws.Cells.Item(115, 14).Value = -43349.332   # N115: -70849.164 -> -43349.332
$ws.Cells.Item(126, 8).Value = 3562.3845   # H126: 3634.25 -> 3562.3845
$ws.Cells.Item(126, 9).Value = 3442.5833   # I126: 3510.0908 -> 3442.5833
$ws.Cells.Item(126, 11).Value = 10327.7499   # K126: 10530.2724 -> 10327.7499
$ws.Cells.Item(126, 13).Value = -7857.749899999999   # M126: -8060.2724 -> -7857.749899999999
$ws.Cells.Item(132, 8).Value = 5565.1816   # H132: 5734.1113 -> 5565.1816
$ws.Cells.Item(132, 9).Value = 5468.5557   # I132: 5658.143 -> 5468.5557
$ws.Cells.Item(132, 11).Value = 16405.6671   # K132: 16974.429 -> 16405.6671
$ws.Cells.Item(132, 13).Value = -13875.6671   # M132: -14444.429 -> -13875.6671
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 2538.2354   # H9: 2812.3333 -> 2538.2354
$ws.Cells.Item(9, 9).Value = 558.125   # I9: 583.3333 -> 558.125
$ws.Cells.Item(9, 11).Value = 1674.375   # K9: 1749.9999 -> 1674.375
$ws.Cells.Item(9, 13).Value = -1450.375   # M9: -1525.9999 -> -1450.375
$ws.Cells.Item(58, 8).Value = 1950   # H58: 1000 -> 1950
$ws.Cells.Item(58, 10).Value = 2900   # J58: 0 -> 2900
$ws.Cells.Item(58, 12).Value = 8700   # L58: 0 -> 8700
$ws.Cells.Item(58, 14).Value = -8956   # N58: None -> -8956
$ws.Cells.Item(113, 8).Value = 378.68182   # H113: 386.22726 -> 378.68182
$ws.Cells.Item(113, 9).Value = 385   # I113: 500 -> 385
$ws.Cells.Item(113, 10).Value = 375.73334   # J113: 352.7647 -> 375.73334
$ws.Cells.Item(113, 11).Value = 1155   # K113: 1500 -> 1155
$ws.Cells.Item(113, 12).Value = 1127.20002   # L113: 1058.2941 -> 1127.20002
$ws.Cells.Item(113, 13).Value = 1015   # M113: 670 -> 1015
$ws.Cells.Item(113, 14).Value = -5467.20002   # N113: -5398.2941 -> -5467.20002
$ws.Cells.Item(131, 8).Value = 49528.24   # H131: 47367.547 -> 49528.24
$ws.Cells.Item(131, 10).Value = 2004.65   # J131: 2004.0952 -> 2004.65
$ws.Cells.Item(131, 12).Value = 6013.950000000001   # L131: 6012.2856 -> 6013.950000000001
$ws.Cells.Item(131, 14).Value = -16093.95   # N131: -16092.2856 -> -16093.95
$ws.Cells.Item(137, 8).Value = 2703.125   # H137: 2847.111 -> 2703.125
$ws.Cells.Item(137, 10).Value = 7000   # J137: 5499.5 -> 7000
$ws.Cells.Item(137, 12).Value = 21000   # L137: 16498.5 -> 21000
$ws.Cells.Item(137, 14).Value = -31200   # N137: -26698.5 -> -31200
$ws.Cells.Item(140, 8).Value = 4545.625   # H140: 4573.625 -> 4545.625
$ws.Cells.Item(140, 9).Value = 4545.625   # I140: 4941.4287 -> 4545.625
$ws.Cells.Item(140, 10).Value = 0   # J140: 1999 -> 0
$ws.Cells.Item(140, 11).Value = 13636.875   # K140: 14824.2861 -> 13636.875
$ws.Cells.Item(140, 12).Value = 0   # L140: 5997 -> 0
$ws.Cells.Item(140, 13).Value = -8456.875   # M140: -9644.286100000001 -> -8456.875
$ws.Cells.Item(140, 14).ClearContents()   # N140: remove (was -16357)
# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7080   # H80: 5050 -> 7080
$ws.Cells.Item(80, 9).Value = 3250   # I80: 2500 -> 3250
$ws.Cells.Item(80, 10).Value = 9633.333000000001   # J80: 7600 -> 9633.333000000001
$ws.Cells.Item(80, 11).Value = 3250   # K80: 2500 -> 3250
$ws.Cells.Item(80, 12).Value = 9633.333000000001   # L80: 7600 -> 9633.333000000001
$ws.Cells.Item(80, 13).Value = -2252   # M80: -1502 -> -2252
$ws.Cells.Item(80, 14).Value = -11629.333   # N80: -9596 -> -11629.333
$ws.Cells.Item(83, 8).Value = 7080   # H83: 5050 -> 7080
$ws.Cells.Item(83, 9).Value = 3250   # I83: 2500 -> 3250
$ws.Cells.Item(83, 10).Value = 9633.333000000001   # J83: 7600 -> 9633.333000000001
$ws.Cells.Item(83, 11).Value = 16250   # K83: 12500 -> 16250
$ws.Cells.Item(83, 12).Value = 48166.665   # L83: 38000 -> 48166.665
$ws.Cells.Item(83, 13).Value = -11258   # M83: -7508 -> -11258
$ws.Cells.Item(83, 14).Value = -58150.665   # N83: -47984 -> -58150.665
$ws.Cells.Item(123, 8).Value = 35674.832   # H123: 35700 -> 35674.832
$ws.Cells.Item(123, 10).Value = 35674.832   # J123: 35700 -> 35674.832
$ws.Cells.Item(123, 12).Value = 35674.832   # L123: 35700 -> 35674.832
$ws.Cells.Item(123, 14).Value = -40574.832   # N123: -40600 -> -40574.832
$ws.Cells.Item(132, 8).Value = 265399.25   # H132: 288065.03 -> 265399.25
$ws.Cells.Item(132, 9).Value = 372591.72   # I132: 386887.53 -> 372591.72
$ws.Cells.Item(132, 10).Value = 2290.4546   # J132: 2577.7778 -> 2290.4546
$ws.Cells.Item(132, 11).Value = 1117775.16   # K132: 1160662.59 -> 1117775.16
$ws.Cells.Item(132, 12).Value = 6871.3638   # L132: 7733.3334 -> 6871.3638
$ws.Cells.Item(132, 13).Value = -1115245.16   # M132: -1158132.59 -> -1115245.16
$ws.Cells.Item(132, 14).Value = -11931.3638   # N132: -12793.3334 -> -11931.3638
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48, 8).Value = 17010.334   # H48: 17360.334 -> 17010.334
$ws.Cells.Item(48, 9).Value = 17010.334   # I48: 17360.334 -> 17010.334
$ws.Cells.Item(48, 11).Value = 17010.334   # K48: 17360.334 -> 17010.334
$ws.Cells.Item(48, 13).Value = -16349.334   # M48: -16699.334 -> -16349.334
$ws.Cells.Item(55, 8).Value = 306.2   # H55: 322.78946 -> 306.2
$ws.Cells.Item(55, 9).Value = 359.7   # I55: 385.1111 -> 359.7
$ws.Cells.Item(55, 10).Value = 252.7   # J55: 266.7 -> 252.7
$ws.Cells.Item(55, 11).Value = 359.7   # K55: 385.1111 -> 359.7
$ws.Cells.Item(55, 12).Value = 252.7   # L55: 266.7 -> 252.7
$ws.Cells.Item(55, 13).Value = -186.7   # M55: -212.1111 -> -186.7
$ws.Cells.Item(55, 14).Value = -598.7   # N55: -612.7 -> -598.7
$ws.Cells.Item(93, 8).Value = 3533.0557   # H93: 3110.9524 -> 3533.0557
$ws.Cells.Item(93, 9).Value = 3592   # I93: 2879.4443 -> 3592
$ws.Cells.Item(93, 10).Value = 3379.8   # J93: 4500 -> 3379.8
$ws.Cells.Item(93, 11).Value = 3592   # K93: 2879.4443 -> 3592
$ws.Cells.Item(93, 12).Value = 3379.8   # L93: 4500 -> 3379.8
$ws.Cells.Item(93, 13).Value = -2344   # M93: -1631.4443 -> -2344
$ws.Cells.Item(93, 14).Value = -5875.8   # N93: -6996 -> -5875.8
$ws.Cells.Item(122, 8).Value = 4276.4546   # H122: 4404.1 -> 4276.4546
$ws.Cells.Item(122, 9).Value = 3630.5   # I122: 3739.2 -> 3630.5
$ws.Cells.Item(122, 10).Value = 5999   # J122: 6398.8 -> 5999
$ws.Cells.Item(122, 11).Value = 10891.5   # K122: 11217.6 -> 10891.5
$ws.Cells.Item(122, 12).Value = 17997   # L122: 19196.4 -> 17997
$ws.Cells.Item(122, 13).Value = -8441.5   # M122: -8767.599999999999 -> -8441.5
$ws.Cells.Item(122, 14).Value = -22897   # N122: -24096.4 -> -22897
$ws.Cells.Item(140, 8).Value = 49983.332   # H140: 45000 -> 49983.332
$ws.Cells.Item(140, 10).Value = 59975   # J140: 60000 -> 59975
$ws.Cells.Item(140, 12).Value = 59975   # L140: 60000 -> 59975
$ws.Cells.Item(140, 14).Value = -70335   # N140: -70360 -> -70335
# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 0   # H11: 1000000 -> 0
$ws.Cells.Item(11, 9).Value = 0   # I11: 1000000 -> 0
$ws.Cells.Item(11, 11).Value = 0   # K11: 1000000 -> 0
$ws.Cells.Item(11, 13).ClearContents()   # M11: remove (was -999858)
$ws.Cells.Item(34, 8).Value = 0   # H34: 10000 -> 0
$ws.Cells.Item(34, 10).Value = 0   # J34: 10000 -> 0
$ws.Cells.Item(34, 12).Value = 0   # L34: 10000 -> 0
$ws.Cells.Item(34, 14).ClearContents()   # N34: remove (was -10406)
$ws.Cells.Item(62, 8).Value = 18282.867   # H62: 17358.938 -> 18282.867
$ws.Cells.Item(62, 10).Value = 10139.8   # J62: 9033.166999999999 -> 10139.8
$ws.Cells.Item(62, 12).Value = 10139.8   # L62: 9033.166999999999 -> 10139.8
$ws.Cells.Item(62, 14).Value = -11387.8   # N62: -10281.167 -> -11387.8
$ws.Cells.Item(65, 8).Value = 18282.867   # H65: 17358.938 -> 18282.867
$ws.Cells.Item(65, 10).Value = 10139.8   # J65: 9033.166999999999 -> 10139.8
$ws.Cells.Item(65, 12).Value = 50699   # L65: 45165.835 -> 50699
$ws.Cells.Item(65, 14).Value = -56939   # N65: -51405.835 -> -56939
$ws.Cells.Item(96, 8).Value = 2621.7144   # H96: 2519 -> 2621.7144
$ws.Cells.Item(96, 9).Value = 1870.4   # I96: 1858.6666 -> 1870.4
$ws.Cells.Item(96, 11).Value = 1870.4   # K96: 1858.6666 -> 1870.4
$ws.Cells.Item(96, 13).Value = -497.4000000000001   # M96: -485.6666 -> -497.4000000000001
$ws.Cells.Item(132, 8).Value = 2476.3225   # H132: 2615.3667 -> 2476.3225
$ws.Cells.Item(132, 9).Value = 2417.4167   # I132: 2596.2173 -> 2417.4167
$ws.Cells.Item(132, 11).Value = 7252.250100000001   # K132: 7788.651899999999 -> 7252.250100000001
$ws.Cells.Item(132, 13).Value = -4722.250100000001   # M132: -5258.651899999999 -> -4722.250100000001
